$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reset the populated region to a clean slate (values + formatting) ---
$ws.Range("A1:H29").Clear()

# --- Rewrite every row of the standardized balance sheet per the new layout ---
# Row 1
$ws.Range("A1").Value = "standardized_balancesheet_label"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").VerticalAlignment = -4160
$ws.Range("B1").Value = "df_Facts_label"
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").VerticalAlignment = -4160
$ws.Range("C1").Value = "df_Fact_Description"
$ws.Range("C1").Font.Bold = $true
$ws.Range("C1").VerticalAlignment = -4160
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").VerticalAlignment = -4160
$ws.Range("E1").Font.Bold = $true
$ws.Range("H1").Font.Bold = $true

# Row 2
$ws.Range("A2").Value = "Cash & Cash Equivalent"
$ws.Range("A2").VerticalAlignment = -4160
$ws.Range("A2").WrapText = $true
$ws.Range("B2").Value = "Cash and Cash Equivalents"
$ws.Range("B2").VerticalAlignment = -4160
$ws.Range("B2").WrapText = $true
$ws.Range("C2").Value = "Amount of currency on hand as well as demand deposits with banks or financial institutions. Includes other kinds of accounts that have the general characteristics of demand deposits. Also includes short-term, highly liquid investments that are both readily convertible to known amounts of cash and so near their maturity that they present insignificant risk of changes in value because of changes in interest rates. Excludes cash and cash equivalents within disposal group and discontinued operation."
$ws.Range("C2").VerticalAlignment = -4160
$ws.Range("C2").WrapText = $true

# Row 3
$ws.Range("A3").Value = "Cash & Cash Equivalent"
$ws.Range("A3").VerticalAlignment = -4160
$ws.Range("A3").WrapText = $true
$ws.Range("B3").Value = "Cash Cash Equivalents Restricted Cash And Restricted Cash Equivalents"
$ws.Range("B3").VerticalAlignment = -4160
$ws.Range("B3").WrapText = $true
$ws.Range("C3").Value = "Amount of cash and cash equivalents, and cash and cash equivalents restricted to withdrawal or usage. Excludes amount for disposal group and discontinued operations. Cash includes, but is not limited to, currency on hand, demand deposits with banks or financial institutions, and other accounts with general characteristics of demand deposits. Cash equivalents include, but are not limited to, short-term, highly liquid investments that are both readily convertible to known amounts of cash and so near their maturity that they present insignificant risk of changes in value because of changes in interest rates."
$ws.Range("C3").VerticalAlignment = -4160
$ws.Range("C3").WrapText = $true

# Row 4
$ws.Range("A4").Value = "Marketable Securities, Current"
$ws.Range("A4").VerticalAlignment = -4160
$ws.Range("A4").WrapText = $true
$ws.Range("B4").Value = "Marketable Securities, Current"
$ws.Range("B4").VerticalAlignment = -4160
$ws.Range("B4").WrapText = $true
$ws.Range("C4").Value = "Amount of investment in marketable security, classified as current."
$ws.Range("C4").VerticalAlignment = -4160
$ws.Range("C4").WrapText = $true

# Row 5
$ws.Range("A5").Value = "Marketable Securities Current"
$ws.Range("A5").VerticalAlignment = -4160
$ws.Range("A5").WrapText = $true
$ws.Range("B5").Value = "Marketable Securities, Current"
$ws.Range("B5").VerticalAlignment = -4160
$ws.Range("B5").WrapText = $true
$ws.Range("C5").Value = "Marketable Securities, Current Amount of investment in marketable security, classified as current."
$ws.Range("C5").VerticalAlignment = -4160
$ws.Range("C5").WrapText = $true

# Row 6
$ws.Range("A6").Value = "Total Accounts Receivable"
$ws.Range("A6").VerticalAlignment = -4160
$ws.Range("A6").WrapText = $true
$ws.Range("B6").Value = "Accounts Receivable, after Allowance for Credit Loss, Current"
$ws.Range("B6").VerticalAlignment = -4160
$ws.Range("B6").WrapText = $true
$ws.Range("C6").Value = "Amount, after allowance for credit loss, of right to consideration from customer for product sold and service rendered in normal course of business, classified as current."
$ws.Range("C6").VerticalAlignment = -4160
$ws.Range("C6").WrapText = $true

# Row 7
$ws.Range("A7").Value = "Total Inventory"
$ws.Range("A7").VerticalAlignment = -4160
$ws.Range("A7").WrapText = $true
$ws.Range("B7").Value = "Inventory, Net"
$ws.Range("B7").VerticalAlignment = -4160
$ws.Range("B7").WrapText = $true
$ws.Range("C7").Value = "Amount after valuation and LIFO reserves of inventory expected to be sold, or consumed within one year or operating cycle, if longer."
$ws.Range("C7").VerticalAlignment = -4160
$ws.Range("C7").WrapText = $true

# Row 8
$ws.Range("A8").Value = "Prepaid Expenses"
$ws.Range("A8").VerticalAlignment = -4160
$ws.Range("A8").WrapText = $true
$ws.Range("B8").Value = "Prepaid Expense, Current"
$ws.Range("B8").VerticalAlignment = -4160
$ws.Range("B8").WrapText = $true
$ws.Range("C8").Value = "Amount of asset related to consideration paid in advance for costs that provide economic benefits within a future period of one year or the normal operating cycle, if longer."
$ws.Range("C8").VerticalAlignment = -4160
$ws.Range("C8").WrapText = $true

# Row 9
$ws.Range("A9").Value = "Other Current Assets"
$ws.Range("A9").VerticalAlignment = -4160
$ws.Range("A9").WrapText = $true
$ws.Range("B9").Value = "Other Assets, Current"
$ws.Range("B9").VerticalAlignment = -4160
$ws.Range("B9").WrapText = $true
$ws.Range("C9").Value = "Other Assets, Current Amount of current assets classified as other."
$ws.Range("C9").VerticalAlignment = -4160
$ws.Range("C9").WrapText = $true

# Row 10
$ws.Range("A10").Value = "Total Current Assets"
$ws.Range("A10").Font.Bold = $true
$ws.Range("A10").VerticalAlignment = -4160
$ws.Range("B10").Value = "Assets, Current"
$ws.Range("B10").Font.Bold = $true
$ws.Range("B10").VerticalAlignment = -4160
$ws.Range("B10").WrapText = $true
$ws.Range("B10").Interior.Color = 65535
$ws.Range("C10").Value = "Sum of the carrying amounts as of the balance sheet date of all assets that are expected to be realized in cash, sold, or consumed within one year (or the normal operating cycle, if longer). Assets are probable future economic benefits obtained or controlled by an entity as a result of past transactions or events."
$ws.Range("C10").VerticalAlignment = -4160
$ws.Range("C10").WrapText = $true

# Row 11
$ws.Range("A11").Value = "Marketable Securities Non Current"
$ws.Range("A11").VerticalAlignment = -4160
$ws.Range("A11").WrapText = $true
$ws.Range("B11").Value = "Marketable Securities, Non Current"
$ws.Range("B11").VerticalAlignment = -4160
$ws.Range("B11").WrapText = $true
$ws.Range("C11").Value = "Amount of investment in marketable security, classified as Non Current."
$ws.Range("C11").VerticalAlignment = -4160
$ws.Range("C11").WrapText = $true

# Row 12
$ws.Range("A12").Value = "Marketable Securities Non Current"
$ws.Range("A12").VerticalAlignment = -4160
$ws.Range("A12").WrapText = $true
$ws.Range("B12").Value = "Marketable Securities, Non Current"
$ws.Range("B12").VerticalAlignment = -4160
$ws.Range("B12").WrapText = $true
$ws.Range("C12").Value = "Marketable Securities, Current Amount of investment in marketable security, classified as non current."
$ws.Range("C12").VerticalAlignment = -4160
$ws.Range("C12").WrapText = $true

# Row 13
$ws.Range("A13").Value = "Property Plant and Equipment"
$ws.Range("A13").VerticalAlignment = -4160
$ws.Range("A13").WrapText = $true
$ws.Range("B13").Value = "Property, Plant and Equipment, Net"
$ws.Range("B13").VerticalAlignment = -4160
$ws.Range("B13").WrapText = $true
$ws.Range("C13").Value = "Amount after accumulated depreciation, depletion and amortization of physical assets used in the normal conduct of business to produce goods and services and not intended for resale. Examples include, but are not limited to, land, buildings, machinery and equipment, office equipment, and furniture and fixtures."
$ws.Range("C13").VerticalAlignment = -4160
$ws.Range("C13").WrapText = $true

# Row 14
$ws.Range("A14").VerticalAlignment = -4160
$ws.Range("A14").WrapText = $true
$ws.Range("B14").Value = "Operating Lease, Right-of-Use Asset"
$ws.Range("B14").VerticalAlignment = -4160
$ws.Range("B14").WrapText = $true
$ws.Range("C14").VerticalAlignment = -4160
$ws.Range("C14").WrapText = $true

# Row 15
$ws.Range("A15").Value = "Property Plant and Equipment"
$ws.Range("A15").VerticalAlignment = -4160
$ws.Range("A15").WrapText = $true
$ws.Range("B15").Value = "Molds and tooling, Net"
$ws.Range("B15").VerticalAlignment = -4160
$ws.Range("B15").WrapText = $true
$ws.Range("C15").Value = "Amount after accumulated depreciation of tangible personal property used to produce goods and services, including, but is not limited to, tools, dies and molds, computer and office equipment."
$ws.Range("C15").VerticalAlignment = -4160
$ws.Range("C15").WrapText = $true

# Row 16
$ws.Range("A16").Value = "Intangible Assets (excl. goodwill)"
$ws.Range("A16").VerticalAlignment = -4160
$ws.Range("A16").WrapText = $true
$ws.Range("B16").Value = "Intangible Assets, Net (Excluding Goodwill)"
$ws.Range("B16").VerticalAlignment = -4160
$ws.Range("B16").WrapText = $true
$ws.Range("C16").Value = "Sum of the carrying amounts of all intangible assets, excluding goodwill, as of the balance sheet date, net of accumulated amortization and impairment charges."
$ws.Range("C16").VerticalAlignment = -4160
$ws.Range("C16").WrapText = $true

# Row 17
$ws.Range("A17").Value = "Goodwill"
$ws.Range("A17").VerticalAlignment = -4160
$ws.Range("B17").Value = "Goodwill"
$ws.Range("B17").VerticalAlignment = -4160
$ws.Range("B17").WrapText = $true
$ws.Range("C17").Value = "Amount after accumulated impairment loss of an asset representing future economic benefits arising from other assets acquired in a business combination that are not individually identified and separately recognized."
$ws.Range("C17").VerticalAlignment = -4160
$ws.Range("C17").WrapText = $true

# Row 18
$ws.Range("A18").Value = "Total Non Current Assets"
$ws.Range("A18").Font.Bold = $true
$ws.Range("A18").VerticalAlignment = -4160
$ws.Range("B18").Value = "Assets, Non Current"
$ws.Range("B18").Font.Bold = $true
$ws.Range("B18").VerticalAlignment = -4160
$ws.Range("B18").WrapText = $true
$ws.Range("C18").Value = "Sum of the carrying amounts as of the balance sheet date of all assets that are expected to be realized in cash, sold or consumed after one year or beyond the normal operating cycle, if longer."
$ws.Range("C18").VerticalAlignment = -4160
$ws.Range("C18").WrapText = $true
$ws.Range("D18").Font.Bold = $true
$ws.Range("D18").VerticalAlignment = -4160

# Row 19
$ws.Range("A19").Value = "Total Assets"
$ws.Range("A19").Font.Bold = $true
$ws.Range("A19").VerticalAlignment = -4160
$ws.Range("B19").Value = "Assets"
$ws.Range("B19").Font.Bold = $true
$ws.Range("B19").VerticalAlignment = -4160
$ws.Range("B19").WrapText = $true
$ws.Range("C19").Value = "Sum of the carrying amounts as of the balance sheet date of all assets that are recognized. Assets are probable future economic benefits obtained or controlled by an entity as a result of past transactions or events."
$ws.Range("C19").VerticalAlignment = -4160
$ws.Range("C19").WrapText = $true
$ws.Range("D19").Font.Bold = $true
$ws.Range("D19").VerticalAlignment = -4160

# Row 20
$ws.Range("A20").Value = "Accounts Payable"
$ws.Range("A20").VerticalAlignment = -4160
$ws.Range("A20").WrapText = $true
$ws.Range("B20").Value = "Accounts Payable, Current"
$ws.Range("B20").VerticalAlignment = -4160
$ws.Range("B20").WrapText = $true
$ws.Range("C20").Value = "Carrying value as of the balance sheet date of liabilities incurred (and for which invoices have typically been received) and payable to vendors for goods and services received that are used in an entity's business. Used to reflect the current portion of the liabilities (due within one year or within the normal operating cycle if longer)."
$ws.Range("C20").VerticalAlignment = -4160
$ws.Range("C20").WrapText = $true

# Row 21
$ws.Range("A21").Value = "Tax Payable"
$ws.Range("A21").VerticalAlignment = -4160
$ws.Range("A21").WrapText = $true
$ws.Range("B21").Value = "Taxes Payable, Current"
$ws.Range("B21").VerticalAlignment = -4160
$ws.Range("B21").WrapText = $true
$ws.Range("C21").Value = "Carrying value as of the balance sheet date of obligations incurred and payable for statutory income, sales, use, payroll, excise, real, property and other taxes. Used to reflect the current portion of the liabilities (due within one year or within the normal operating cycle if longer)."
$ws.Range("C21").VerticalAlignment = -4160
$ws.Range("C21").WrapText = $true

# Row 22
$ws.Range("A22").Value = "Short Term Debt"
$ws.Range("A22").VerticalAlignment = -4160
$ws.Range("A22").WrapText = $true
$ws.Range("B22").Value = "Long Term Debt, Current Maturities"
$ws.Range("B22").VerticalAlignment = -4160
$ws.Range("B22").WrapText = $true
$ws.Range("C22").Value = "Amount, after unamortized (discount) premium and debt issuance costs, of Long Term debt, classified as current. Includes, but not limited to, notes payable, bonds payable, debentures, mortgage loans and commercial paper. Excludes capital lease obligations."
$ws.Range("C22").VerticalAlignment = -4160
$ws.Range("C22").WrapText = $true

# Row 23
$ws.Range("A23").Value = "Operating Lease Liability Current"
$ws.Range("A23").VerticalAlignment = -4160
$ws.Range("A23").WrapText = $true
$ws.Range("B23").Value = "Operating Lease, Liability, Current"
$ws.Range("B23").VerticalAlignment = -4160
$ws.Range("B23").WrapText = $true
$ws.Range("C23").Value = "Present value of lessee's discounted obligation for lease payments from operating lease, classified as current."
$ws.Range("C23").VerticalAlignment = -4160
$ws.Range("C23").WrapText = $true

# Row 24
$ws.Range("A24").Value = "Finance Lease Liability Current"
$ws.Range("A24").VerticalAlignment = -4160
$ws.Range("A24").WrapText = $true
$ws.Range("B24").Value = "Finance Lease, Liability, Current"
$ws.Range("B24").VerticalAlignment = -4160
$ws.Range("B24").WrapText = $true
$ws.Range("C24").Value = "Present value of lessees discounted obligation for lease payments from finance lease, classified as current."
$ws.Range("C24").VerticalAlignment = -4160
$ws.Range("C24").WrapText = $true

# Row 25
$ws.Range("A25").Value = "Total Current Liabilities"
$ws.Range("A25").Font.Bold = $true
$ws.Range("A25").VerticalAlignment = -4160
$ws.Range("B25").Value = "Liabilities, Current"
$ws.Range("B25").Font.Bold = $true
$ws.Range("B25").VerticalAlignment = -4160
$ws.Range("B25").WrapText = $true
$ws.Range("C25").Value = "Total obligations incurred as part of normal operations that are expected to be paid during the following twelve months or within one business cycle, if longer."
$ws.Range("C25").VerticalAlignment = -4160
$ws.Range("C25").WrapText = $true
$ws.Range("D25").Font.Bold = $true
$ws.Range("D25").VerticalAlignment = -4160

# Row 26
$ws.Range("A26").Value = "Non Current Debts"
$ws.Range("A26").VerticalAlignment = -4160
$ws.Range("A26").WrapText = $true
$ws.Range("B26").Value = "Long Term Debt, Excluding Current Maturities"
$ws.Range("B26").VerticalAlignment = -4160
$ws.Range("B26").WrapText = $true
$ws.Range("C26").Value = "Amount after unamortized (discount) premium and debt issuance costs of Long Term debt classified as Non Current and excluding amounts to be repaid within one year or the normal operating cycle, if longer. Includes, but not limited to, notes payable, bonds payable, debentures, mortgage loans and commercial paper. Excludes capital lease obligation."
$ws.Range("C26").VerticalAlignment = -4160
$ws.Range("C26").WrapText = $true

# Row 27
$ws.Range("A27").Value = "Operating Lease Liability Non Current"
$ws.Range("A27").VerticalAlignment = -4160
$ws.Range("A27").WrapText = $true
$ws.Range("B27").Value = "Operating Lease, Liability, Non Current"
$ws.Range("B27").VerticalAlignment = -4160
$ws.Range("B27").WrapText = $true
$ws.Range("C27").Value = "Present value of lessees discounted obligation for lease payments from operating lease, classified as Non Current."
$ws.Range("C27").VerticalAlignment = -4160
$ws.Range("C27").WrapText = $true

# Row 28
$ws.Range("A28").Value = "Finance Lease Liability Non Current"
$ws.Range("A28").VerticalAlignment = -4160
$ws.Range("A28").WrapText = $true
$ws.Range("B28").Value = "Finance Lease, Liability, Current"
$ws.Range("B28").VerticalAlignment = -4160
$ws.Range("B28").WrapText = $true
$ws.Range("C28").Value = "Present value of lessees discounted obligation for lease payments from finance lease, classified as Non Current."
$ws.Range("C28").VerticalAlignment = -4160
$ws.Range("C28").WrapText = $true

# Row 29
$ws.Range("A29").Value = "Total Non Current Liabilities"
$ws.Range("A29").Font.Bold = $true
$ws.Range("A29").VerticalAlignment = -4160
$ws.Range("B29").Value = "Liabilities, Non Current"
$ws.Range("B29").Font.Bold = $true
$ws.Range("B29").VerticalAlignment = -4160
$ws.Range("B29").WrapText = $true
$ws.Range("C29").Value = "Amount of obligation due after one year or beyond the normal operating cycle, if longer."
$ws.Range("C29").VerticalAlignment = -4160
$ws.Range("C29").WrapText = $true
$ws.Range("D29").Font.Bold = $true
$ws.Range("D29").VerticalAlignment = -4160

# Row 30
$ws.Range("A30").Value = "Total Liabilities"
$ws.Range("A30").Font.Bold = $true
$ws.Range("A30").VerticalAlignment = -4160
$ws.Range("B30").Value = "Liabilities"
$ws.Range("B30").Font.Bold = $true
$ws.Range("B30").VerticalAlignment = -4160
$ws.Range("B30").WrapText = $true
$ws.Range("C30").Value = "Sum of the carrying amounts as of the balance sheet date of all liabilities that are recognized. Liabilities are probable future sacrifices of economic benefits arising from present obligations of an entity to transfer assets or provide services to other entities in the future."
$ws.Range("C30").VerticalAlignment = -4160
$ws.Range("C30").WrapText = $true
$ws.Range("D30").Font.Bold = $true
$ws.Range("D30").VerticalAlignment = -4160

# Row 31
$ws.Range("A31").Value = "Preferred Stock"
$ws.Range("A31").VerticalAlignment = -4160
$ws.Range("A31").WrapText = $true
$ws.Range("B31").Value = "Preferred Stock, Shares Outstanding"
$ws.Range("B31").VerticalAlignment = -4160
$ws.Range("B31").WrapText = $true
$ws.Range("C31").Value = "Aggregate share number for all nonredeemable preferred stock (or preferred stock redeemable solely at the option of the issuer) held by stockholders. Does not include preferred shares that have been repurchased."
$ws.Range("C31").VerticalAlignment = -4160
$ws.Range("C31").WrapText = $true

# Row 32
$ws.Range("A32").Value = "Retained Earnings"
$ws.Range("A32").VerticalAlignment = -4160
$ws.Range("A32").WrapText = $true
$ws.Range("B32").Value = "Retained Earnings (Accumulated Deficit)"
$ws.Range("B32").VerticalAlignment = -4160
$ws.Range("B32").WrapText = $true
$ws.Range("C32").Value = "The cumulative amount of the reporting entity's undistributed earnings or deficit."
$ws.Range("C32").VerticalAlignment = -4160
$ws.Range("C32").WrapText = $true

# Row 33
$ws.Range("A33").Value = "Accumulated other comprehensive income (loss)"
$ws.Range("A33").VerticalAlignment = -4160
$ws.Range("A33").WrapText = $true
$ws.Range("B33").Value = "Accumulated Other Comprehensive Income (Loss), Net of Tax"
$ws.Range("B33").VerticalAlignment = -4160
$ws.Range("B33").WrapText = $true
$ws.Range("C33").Value = "Accumulated change in equity from transactions and other events and circumstances from non-owner sources, net of tax effect, at period end. Excludes Net Income (Loss), and accumulated changes in equity from transactions resulting from investments by owners and distributions to owners. Includes foreign currency translation items, certain pension adjustments, unrealized gains and losses on certain investments in debt and equity securities, other than temporary impairment (OTTI) losses related to factors other than credit losses on available-for-sale and held-to-maturity debt securities that an entity does not intend to sell and it is not more likely than not that the entity will be required to sell before recovery of the amortized cost basis, as well as changes in the fair value of derivatives related to the effective portion of a designated cash flow hedge."
$ws.Range("C33").VerticalAlignment = -4160
$ws.Range("C33").WrapText = $true

# Row 34
$ws.Range("A34").Value = "Minority interest"
$ws.Range("A34").VerticalAlignment = -4160
$ws.Range("A34").WrapText = $true
$ws.Range("B34").Value = "Stockholders Equity Attributable to Noncontrolling Interest"
$ws.Range("B34").VerticalAlignment = -4160
$ws.Range("B34").WrapText = $true
$ws.Range("C34").Value = "Total of all stockholders equity (deficit) items, net of receivables from officers, directors, owners, and affiliates of the entity which is directly or indirectly attributable to that ownership interest in subsidiary equity which is not attributable to the parent (that is, noncontrolling interest, previously referred to as minority interest)."
$ws.Range("C34").VerticalAlignment = -4160
$ws.Range("C34").WrapText = $true

# Row 35
$ws.Range("A35").Value = "Total Stockholders Equity"
$ws.Range("A35").Font.Bold = $true
$ws.Range("A35").VerticalAlignment = -4160
$ws.Range("B35").Value = "Stockholders Equity, Including Portion Attributable to Noncontrolling Interest"
$ws.Range("B35").Font.Bold = $true
$ws.Range("B35").VerticalAlignment = -4160
$ws.Range("B35").WrapText = $true
$ws.Range("C35").Value = "Amount of stockholders' equity (deficit), net of receivables from officers, directors, owners, and affiliates of the entity, attributable to both the parent and noncontrolling interests. Amount excludes temporary equity. Alternate caption for the concept is permanent equity."
$ws.Range("C35").VerticalAlignment = -4160
$ws.Range("C35").WrapText = $true

# Row 36
$ws.Range("A36").Value = "Total Liabilities & Stockholders Equity"
$ws.Range("A36").Font.Bold = $true
$ws.Range("A36").VerticalAlignment = -4160
$ws.Range("B36").Value = "Liabilities and Equity"
$ws.Range("B36").Font.Bold = $true
$ws.Range("B36").VerticalAlignment = -4160
$ws.Range("C36").Value = "Amount of liabilities and equity items, including the portion of equity attributable to noncontrolling interests, if any."
$ws.Range("C36").VerticalAlignment = -4160
$ws.Range("C36").WrapText = $true
$ws.Range("D36").Font.Bold = $true
$ws.Range("D36").VerticalAlignment = -4160

# --- Row heights ---
$ws.Rows.Item(1).RowHeight = 16
$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(3).RowHeight = 15
$ws.Rows.Item(4).RowHeight = 15
$ws.Rows.Item(5).RowHeight = 15
$ws.Rows.Item(6).RowHeight = 15
$ws.Rows.Item(7).RowHeight = 15
$ws.Rows.Item(8).RowHeight = 15
$ws.Rows.Item(9).RowHeight = 15
$ws.Rows.Item(10).RowHeight = 15
$ws.Rows.Item(11).RowHeight = 15
$ws.Rows.Item(12).RowHeight = 15
$ws.Rows.Item(13).RowHeight = 15
$ws.Rows.Item(14).RowHeight = 15
$ws.Rows.Item(15).RowHeight = 15
$ws.Rows.Item(16).RowHeight = 15
$ws.Rows.Item(17).RowHeight = 15
$ws.Rows.Item(18).RowHeight = 15
$ws.Rows.Item(19).RowHeight = 15
$ws.Rows.Item(20).RowHeight = 15
$ws.Rows.Item(21).RowHeight = 15
$ws.Rows.Item(22).RowHeight = 15
$ws.Rows.Item(23).RowHeight = 15
$ws.Rows.Item(24).RowHeight = 15
$ws.Rows.Item(25).RowHeight = 15
$ws.Rows.Item(26).RowHeight = 15
$ws.Rows.Item(27).RowHeight = 15
$ws.Rows.Item(28).RowHeight = 15
$ws.Rows.Item(29).RowHeight = 15
$ws.Rows.Item(30).RowHeight = 15
$ws.Rows.Item(31).RowHeight = 18
$ws.Rows.Item(32).RowHeight = 15
$ws.Rows.Item(33).RowHeight = 15
$ws.Rows.Item(34).RowHeight = 15
$ws.Rows.Item(35).RowHeight = 15
$ws.Rows.Item(36).RowHeight = 34

# --- View / selection state ---
$ws.Activate()
$aw = $excel.ActiveWindow
$aw.ScrollColumn = 2
$aw.ScrollRow = 1
$ws.Range("B10").Select()

